$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) New paragraph: "Milgram & Kishino (1994) erwähnen Möglichkeiten von
#    Audio, Haptic, Vestibular AR" inserted right after the
#    'Dictionary: "Vision technologies ... real-world scene."' paragraph
#    and before the 'Durlach, Ternier (2012) (Audio VR)' paragraph.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute('real-world scene."', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorIndex = $find.Parent.Paragraphs.Item(1).Index
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newParaRange = $newPara.Range
$newParaRange.Collapse(1)

$milgramXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:pStyle w:val="StandardErstzeileneinzug"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Milgram</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> &amp; </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Kishino</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> (1994) erwähnen Möglichkeiten von Audio</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Haptic</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Vestibular</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> AR</w:t></w:r>' + `
  '<w:bookmarkStart w:id="7" w:name="_GoBack"/><w:bookmarkEnd w:id="7"/>' + `
  '</w:p>'

$newParaRange.InsertXML($milgramXml)

Write-Output "step1 done"
